$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = "somatotipo"
$ws.Range("F15").Value = "dieta"

$ws.Range("E16").Value = "endomorfo"
$ws.Range("F16").Value = "carboidrati"

$ws.Range("E17").Value = "mesomorfo"
$ws.Range("F17").Value = "proteine"

$ws.Range("E18").Value = "ectomorfo"
$ws.Range("F18").Value = "zuccheri"
